$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.621.31'
$ws.Range('E2').Value = '  +5.07%  '

$ws.Range('D3').Value = '3.323.60'
$ws.Range('E3').Value = '  +4.64%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '553.78'
$ws.Range('E5').Value = '  +3.59%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.87'
$ws.Range('E6').Value = '  +4.44%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('E8').Value = '  +1.18%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.48'
$ws.Range('E9').Value = '  +2.76%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.117'
$ws.Range('E10').Value = '  +3.90%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.434'
$ws.Range('E11').Value = '  +1.30%  '

$ws.Range('D12').Value = '3.899.55'
$ws.Range('E12').Value = '  +4.78%  '

$ws.Range('E13').Value = '  -0.85%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000180'
$ws.Range('E14').Value = '  +4.42%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.83'
$ws.Range('E15').Value = '  +3.47%  '

$ws.Range('D16').Value = '62.636.11'
$ws.Range('E16').Value = '  +5.13%  '

$ws.Range('D17').Value = '3.333.15'
$ws.Range('E17').Value = '  +5.31%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.50'
$ws.Range('E18').Value = '  +5.05%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.74'
$ws.Range('E19').Value = '  +5.74%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.46'
$ws.Range('E20').Value = '  +3.08%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '386.81'
$ws.Range('E21').Value = '  +2.04%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.25%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.536'
$ws.Range('E23').Value = '  +1.22%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.90'
$ws.Range('E24').Value = '  +0.95%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.178'
$ws.Range('E25').Value = '  +4.52%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.85'
$ws.Range('E26').Value = '  -0.72%  '

$ws.Range('D27').Value = '0.0₃0964'
$ws.Range('E27').Value = '  +6.90%  '

$ws.Range('E28').Value = '  +0.17%  '

$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.97'
$ws.Range('E29').Value = '  +3.67%  '

$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.43'
$ws.Range('E30').Value = '  +4.59%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.91'
$ws.Range('E31').Value = '  +2.62%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.55'
$ws.Range('E32').Value = '  +4.19%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.30'
$ws.Range('E33').Value = '  +9.80%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.71'
$ws.Range('E34').Value = '  +4.01%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.49'
$ws.Range('E35').Value = '  +10.10%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '160.55'
$ws.Range('E36').Value = '  +2.52%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.88'
$ws.Range('E37').Value = '  +12.17%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '27.02'
$ws.Range('E38').Value = '  +6.41%  '

$ws.Range('D39').Value = '2.856.85'
$ws.Range('E39').Value = '  +4.10%  '

$ws.Range('E40').Value = '  +3.50%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0314'
$ws.Range('E41').Value = '  +8.56%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.31'
$ws.Range('E42').Value = '  +0.90%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.751'
$ws.Range('E43').Value = '  +3.69%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.77'
$ws.Range('E44').Value = '  +3.44%  '

$ws.Range('E45').Value = '  +4.33%  '

$ws.Range('B46').Value = 'RenzoRestakedETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D46').Value = '3.370.34'
$ws.Range('E46').Value = '  +4.70%  '

$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '21.93'
$ws.Range('E47').Value = '  +7.22%  '

$ws.Range('E48').Value = '  +3.30%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.29'
$ws.Range('E49').Value = '  +1.85%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.803'
$ws.Range('E50').Value = '  +4.07%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '282.55'
$ws.Range('E51').Value = '  +8.83%  '
